$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (header "Förändrad") rows 2 through 120 all hold the serial date
# 45185 (2023-09-16) and must be updated to 45204 (2023-10-05).
$ws.Range("C2:C120").Value2 = 45204
